$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the sentence "...characterize the thermodynamic systems at thermal
# equilibrium." used to be split across two runs around a stray "_GoBack"
# bookmark (left over from a previous edit session, splitting the word
# "systems" into "s" + "ystems"). Re-typing the (unchanged) text collapses
# it back into a single run and drops that bookmark.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute(
    "thermodynamic systems at thermal equilibrium.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    # Force a genuine content change first (Word/this host no-ops a
    # same-text assignment), then restore the real text so the run is
    # rebuilt as a single run with no embedded bookmark.
    $rng1.Text = "TEMP_PLACEHOLDER_TEXT"

    $rng1b = $d.Content
    $rng1b.Find.ClearFormatting()
    $rng1b.Find.Execute("TEMP_PLACEHOLDER_TEXT", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng1b.Text = "thermodynamic systems at thermal equilibrium."

    # Side effect of the host re-normalizing everything after the removed
    # bookmark: the next two (unrelated, identically-formatted) runs in the
    # same paragraph - " In all, these laws " and "describe how these
    # quantities..." - get coalesced into one run too. That pairing was not
    # part of the source edit, so split it back apart the same way any new
    # run boundary gets created here: toggle Bold off then back on over the
    # second half only, which leaves formatting untouched but forces a
    # fresh run to start at that position.
    $rngSplit = $d.Content
    $rngSplit.Find.ClearFormatting()
    $foundSplit = $rngSplit.Find.Execute(
        "describe how these quantities would behave under various " + `
        "circumstances and forbid certain phenomena from occurring such " + `
        "as perpetual motion. ",
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundSplit) {
        $tail = $d.Range($rngSplit.Start, $rngSplit.End)
        $origBold = $tail.Bold
        $tail.Bold = 1
        $tail2 = $d.Range($rngSplit.Start, $rngSplit.End)
        $tail2.Bold = $origBold
    }
}

# ---------------------------------------------------------------------------
# Edit 2: "...will be in the following chapter." gains a parenthetical
# " (chapter 3)" right before the final period, and the cursor's last
# position (the "_GoBack" bookmark) moves to sit right after the new text,
# before the period.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute(
    "will be in the following chapter.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $periodPos = $rng2.End - 1
    $insertion = " (chapter 3)"

    $insertRange = $d.Range($periodPos, $periodPos)
    $insertRange.InsertBefore($insertion)

    # Split the newly inserted text into its own run (distinct from the run
    # that precedes it) without altering any visible formatting: toggle
    # Bold off then back on to force the host to re-segment runs.
    $newRunRange = $d.Range($periodPos, $periodPos + $insertion.Length)
    $origBold = $newRunRange.Bold
    $newRunRange.Bold = 1
    $newRunRange2 = $d.Range($periodPos, $periodPos + $insertion.Length)
    $newRunRange2.Bold = $origBold

    # Re-point the "_GoBack" bookmark (last edit position) to just after the
    # inserted text, ahead of the trailing period.
    $bmPos = $periodPos + $insertion.Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
